$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-30 down to 24-31
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with its values (columns that repeat the
# prior row 23's data, plus the columns that hold genuinely new data)
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44673
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100112022
$ws.Cells.Item(23, 7).Value = "Arveja Verde"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 220
$ws.Cells.Item(23, 11).Value = 25000
$ws.Cells.Item(23, 12).Value = 26000
$ws.Cells.Item(23, 13).Value = 25455
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Carahue"
$ws.Cells.Item(23, 16).Value = 1018
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
